$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently contains duplicate rows for "Swine / pigs" (rows 19-20)
# and "Turkeys" (rows 21-22), followed by "Humans" (row 23).
# Remove the duplicate rows (22 then 20), shifting remaining rows up so the
# final layout is: ... row19 Swine/pigs, row20 Turkeys, row21 Humans.

$ws.Rows.Item(22).Delete()
$ws.Rows.Item(20).Delete()
